$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1278.4667
$ws.Range("I40").Value = 723.3333
$ws.Range("J40").Value = 1648.5555
$ws.Range("K40").Value = 723.3333
$ws.Range("L40").Value = 1648.5555
$ws.Range("M40").Value = -548.3333
$ws.Range("N40").Value = -1998.5555

$ws.Range("H74").Value = 6253836
$ws.Range("I74").Value = 3747.3333
$ws.Range("J74").Value = 8932446
$ws.Range("K74").Value = 3747.3333
$ws.Range("L74").Value = 8932446
$ws.Range("M74").Value = -2811.3333
$ws.Range("N74").Value = -8934318

$ws.Range("H77").Value = 6253836
$ws.Range("I77").Value = 3747.3333
$ws.Range("J77").Value = 8932446
$ws.Range("K77").Value = 18736.6665
$ws.Range("L77").Value = 44662230
$ws.Range("M77").Value = -14056.6665
$ws.Range("N77").Value = -44671590

$ws.Range("H137").Value = 1529.9445
$ws.Range("I137").Value = 1524.0834
$ws.Range("K137").Value = 4572.2502
$ws.Range("M137").Value = -2022.2502

$ws.Range("H139").Value = 50120
$ws.Range("J139").Value = 50120
$ws.Range("L139").Value = 50120
$ws.Range("N139").Value = -60400

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 569.6177
$ws.Range("I2").Value = 582.3871
$ws.Range("J2").Value = 437.66666
$ws.Range("K2").Value = 582.3871
$ws.Range("L2").Value = 437.66666
$ws.Range("M2").Value = -469.3871
$ws.Range("N2").Value = -663.66666

$ws.Range("H45").Value = 3561.4583
$ws.Range("I45").Value = 3233.0833
$ws.Range("J45").Value = 3889.8333
$ws.Range("K45").Value = 3233.0833
$ws.Range("L45").Value = 3889.8333
$ws.Range("M45").Value = -2856.0833
$ws.Range("N45").Value = -4643.8333

$ws.Range("H74").Value = 25002460
$ws.Range("I74").Value = 37039130
$ws.Range("J74").Value = 3231.6924
$ws.Range("K74").Value = 37039130
$ws.Range("L74").Value = 3231.6924
$ws.Range("M74").Value = -37038256
$ws.Range("N74").Value = -4979.6924

$ws.Range("H77").Value = 25002460
$ws.Range("I77").Value = 37039130
$ws.Range("J77").Value = 3231.6924
$ws.Range("K77").Value = 185195650
$ws.Range("L77").Value = 16158.462
$ws.Range("M77").Value = -185191282
$ws.Range("N77").Value = -24894.462

$ws.Range("H88").Value = 112806.78
$ws.Range("J88").Value = 144323
$ws.Range("L88").Value = 144323
$ws.Range("N88").Value = -145135

$ws.Range("H91").Value = 112806.78
$ws.Range("J91").Value = 144323
$ws.Range("L91").Value = 144323
$ws.Range("N91").Value = -147131

$ws.Range("H97").Value = 47620224
$ws.Range("I97").Value = 991.13336
$ws.Range("K97").Value = 991.13336
$ws.Range("M97").Value = -495.13336

$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

$ws.Range("H116").Value = 569.6177
$ws.Range("I116").Value = 582.3871
$ws.Range("J116").Value = 437.66666
$ws.Range("K116").Value = 582.3871
$ws.Range("L116").Value = 437.66666
$ws.Range("M116").Value = 1711.6129
$ws.Range("N116").Value = -5025.66666

$ws.Range("H122").Value = 1263
$ws.Range("I122").Value = 1263
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3789
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1339
$ws.Range("N122").ClearContents()

$ws.Range("H138").Value = 50219
$ws.Range("J138").Value = 50219
$ws.Range("L138").Value = 50219
$ws.Range("N138").Value = -60499

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 569.6177
$ws.Range("I3").Value = 582.3871
$ws.Range("J3").Value = 437.66666
$ws.Range("K3").Value = 582.3871
$ws.Range("L3").Value = 437.66666
$ws.Range("M3").Value = -468.3871
$ws.Range("N3").Value = -665.66666

$ws.Range("H25").Value = 442.66666
$ws.Range("I25").Value = 442.66666
$ws.Range("K25").Value = 442.66666
$ws.Range("M25").Value = -207.66666

$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws.Range("H81").Value = 12308.571
$ws.Range("J81").Value = 12308.571
$ws.Range("L81").Value = 12308.571
$ws.Range("N81").Value = -14430.571

$ws.Range("H84").Value = 12308.571
$ws.Range("J84").Value = 12308.571
$ws.Range("L84").Value = 36925.713
$ws.Range("N84").Value = -47533.713

$ws.Range("H88").Value = 20000
$ws.Range("J88").Value = 20000
$ws.Range("L88").Value = 20000
$ws.Range("N88").Value = -20812

$ws.Range("H91").Value = 20000
$ws.Range("J91").Value = 20000
$ws.Range("L91").Value = 20000
$ws.Range("N91").Value = -22808

$ws.Range("H134").Value = 4935.5454
$ws.Range("I134").Value = 4616.893
$ws.Range("J134").Value = 6720
$ws.Range("K134").Value = 13850.679
$ws.Range("L134").Value = 20160
$ws.Range("M134").Value = -11315.679
$ws.Range("N134").Value = -25230

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3532.0334
$ws.Range("I31").Value = 841.625
$ws.Range("J31").Value = 6606.7856
$ws.Range("K31").Value = 841.625
$ws.Range("L31").Value = 6606.7856
$ws.Range("M31").Value = -546.625
$ws.Range("N31").Value = -7196.7856

$ws.Range("H34").Value = 3532.0334
$ws.Range("I34").Value = 841.625
$ws.Range("J34").Value = 6606.7856
$ws.Range("K34").Value = 841.625
$ws.Range("L34").Value = 6606.7856
$ws.Range("M34").Value = -639.625
$ws.Range("N34").Value = -7010.7856

$ws.Range("H134").Value = 1847.5
$ws.Range("I134").Value = 1771.25
$ws.Range("K134").Value = 5313.75
$ws.Range("M134").Value = -2778.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H24").Value = 638
$ws.Range("J24").Value = 774
$ws.Range("L24").Value = 2322
$ws.Range("N24").Value = -2782

$ws.Range("H113").Value = 831.2941
$ws.Range("J113").Value = 1171.6666
$ws.Range("L113").Value = 3514.9998
$ws.Range("N113").Value = -7854.9998

$ws.Range("H131").Value = 704.36
$ws.Range("J131").Value = 730.6087
$ws.Range("L131").Value = 2191.8261
$ws.Range("N131").Value = -12271.8261

$ws.Range("H132").Value = 528.75
$ws.Range("I132").Value = 528.75
$ws.Range("K132").Value = 4758.75
$ws.Range("M132").Value = -2228.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1728.04
$ws.Range("I102").Value = 1378.3182
$ws.Range("J102").Value = 4292.6665
$ws.Range("K102").Value = 1378.3182
$ws.Range("L102").Value = 4292.6665
$ws.Range("M102").Value = 243.6818000000001
$ws.Range("N102").Value = -7536.6665

$ws.Range("H122").Value = 2464.1428
$ws.Range("I122").Value = 1069.8
$ws.Range("K122").Value = 3209.4
$ws.Range("M122").Value = -759.3999999999996

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1165.0714
$ws.Range("I82").Value = 1116.375
$ws.Range("J82").Value = 1230
$ws.Range("K82").Value = 1116.375
$ws.Range("L82").Value = 1230
$ws.Range("M82").Value = -755.375
$ws.Range("N82").Value = -1952

$ws.Range("H85").Value = 1165.0714
$ws.Range("I85").Value = 1116.375
$ws.Range("J85").Value = 1230
$ws.Range("K85").Value = 1116.375
$ws.Range("L85").Value = 1230
$ws.Range("M85").Value = 131.625
$ws.Range("N85").Value = -3726

$ws.Range("H132").Value = 2851.5386
$ws.Range("I132").Value = 1761.5
$ws.Range("K132").Value = 5284.5
$ws.Range("M132").Value = -2754.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2155.3333
$ws.Range("I81").Value = 1514.6666
$ws.Range("J81").Value = 2796
$ws.Range("K81").Value = 3029.3332
$ws.Range("L81").Value = 5592
$ws.Range("M81").Value = -1968.3332
$ws.Range("N81").Value = -7714

$ws.Range("H84").Value = 2155.3333
$ws.Range("I84").Value = 1514.6666
$ws.Range("J84").Value = 2796
$ws.Range("K84").Value = 15146.666
$ws.Range("L84").Value = 27960
$ws.Range("M84").Value = -9842.666000000001
$ws.Range("N84").Value = -38568

$ws.Range("H126").Value = 1473.4117
$ws.Range("I126").Value = 1099.8928
$ws.Range("K126").Value = 3299.6784
$ws.Range("M126").Value = -829.6784000000002

$ws.Range("H139").Value = 50881.332
$ws.Range("J139").Value = 50881.332
$ws.Range("L139").Value = 50881.332
$ws.Range("N139").Value = -61161.332
